$wb = $excel.ActiveWorkbook

# Update both "展览" and "全部类型" sheets with the same new values
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 8928
    $ws.Range("F4").Value = 443
    $ws.Range("F5").Value = 449
}
